$d = $word.ActiveDocument

# The edit trims the document-wide rPrDefault/pPrDefault in styles.xml so
# that only the non-default values survive (everything else in the diff
# simply matches Word's built-in defaults and is dropped).
#
# docDefaults isn't reachable through the normal Style/Font/ParagraphFormat
# object-model surface, so we round-trip the package through WordOpenXML
# (the flat-OPC textual representation of the whole .docx, including
# word/styles.xml) and do a literal string replace on the <w:docDefaults>
# block.

$oldDocDefaults = '<w:docDefaults><w:rPrDefault><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b w:val="0"/><w:i w:val="0"/><w:smallCaps w:val="0"/><w:strike w:val="0"/><w:color w:val="000000"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:u w:val="none"/><w:shd w:val="clear" w:fill="auto"/><w:vertAlign w:val="baseline"/><w:lang w:val="en"/></w:rPr></w:rPrDefault><w:pPrDefault><w:pPr><w:keepNext w:val="0"/><w:keepLines w:val="0"/><w:widowControl/><w:pBdr><w:top w:val="nil" w:sz="0" w:space="0"/><w:left w:val="nil" w:sz="0" w:space="0"/><w:bottom w:val="nil" w:sz="0" w:space="0"/><w:right w:val="nil" w:sz="0" w:space="0"/><w:between w:val="nil" w:sz="0" w:space="0"/></w:pBdr><w:shd w:val="clear" w:fill="auto"/><w:spacing w:before="0" w:after="0" w:line="276" w:lineRule="auto"/><w:ind w:left="0" w:right="0" w:firstLine="0"/><w:contextualSpacing w:val="0"/><w:jc w:val="left"/></w:pPr></w:pPrDefault></w:docDefaults>'

$newDocDefaults = '<w:docDefaults><w:rPrDefault><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="en"/></w:rPr></w:rPrDefault><w:pPrDefault><w:pPr><w:spacing w:line="276" w:lineRule="auto"/></w:pPr></w:pPrDefault></w:docDefaults>'

$xml = $d.WordOpenXML

if ($xml.Contains($oldDocDefaults)) {
    $xml = $xml.Replace($oldDocDefaults, $newDocDefaults)
    $d.WordOpenXML = $xml
    Write-Output "docDefaults simplified"
} else {
    Write-Output "WARNING: expected docDefaults block not found"
}
